# Swap the slide-master theme's 12 color-scheme slots from the
# "Integral" palette over to the stock "Office" palette (the palette
# that the deck's notes-master theme already uses).
#
# Target RGB values (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink),
# taken from the "Office" theme:
#   dk1=000000 lt1=FFFFFF dk2=44546A lt2=E7E6E6
#   accent1=5B9BD5 accent2=ED7D31 accent3=A5A5A5 accent4=FFC000
#   accent5=4472C4 accent6=70AD47 hlink=0563C1 folHlink=954F72

$p = $ppt.ActivePresentation
$tcs = $p.SlideMaster.Theme.ThemeColorScheme

function Set-ThemeRGB($index, $r, $g, $b) {
    $tcs.Colors($index).RGB = $r + ($g * 256) + ($b * 65536)
}

Set-ThemeRGB 1  0   0   0     # dk1
Set-ThemeRGB 2  255 255 255   # lt1
Set-ThemeRGB 3  68  84  106   # dk2
Set-ThemeRGB 4  231 230 230   # lt2
Set-ThemeRGB 5  91  155 213   # accent1
Set-ThemeRGB 6  237 125 49    # accent2
Set-ThemeRGB 7  165 165 165   # accent3
Set-ThemeRGB 8  255 192 0     # accent4
Set-ThemeRGB 9  68  114 196   # accent5
Set-ThemeRGB 10 112 173 71    # accent6
Set-ThemeRGB 11 5   99  193   # hlink
Set-ThemeRGB 12 149 79  114   # folHlink
